# Re-run experiments/2025-09-05 since previous run had spectral norm loss
# with lambda 1e-8 -> append the new training-config log entry, and make
# sure every existing data row in the "training" sheet has a (blank) value
# in the "note" column, matching the rest of the log sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("training")

# The "note" column (D) is sparse/blank for every existing run, but a real
# (empty) cell should be present for each logged row - backfill it for the
# current rows (2-20) the same way it already exists on the other log
# sheets.
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 4).Value = "'"
}

# Append the new log row (21) for the re-run training config.
$row = 21
$ws.Cells.Item($row, 1).Value  = "2025-09-05 15:27:52"
$ws.Cells.Item($row, 2).Value  = "training"
$ws.Cells.Item($row, 3).Value  = "configs/training/2025-09-05/a/0001"
$ws.Cells.Item($row, 4).Value  = "'"
$ws.Cells.Item($row, 5).Value  = "['cross_entropy']"
$ws.Cells.Item($row, 6).Value  = "[1.0]"
$ws.Cells.Item($row, 7).Value  = "['torch.optim.adamw.AdamW']"
$ws.Cells.Item($row, 8).Value  = "[0.001]"
$ws.Cells.Item($row, 9).Value  = 128
$ws.Cells.Item($row, 10).Value = 128
$ws.Cells.Item($row, 11).Value = "general_utils.ml.training.NoImprovementStopping"
$ws.Cells.Item($row, 12).Value = 8
$ws.Cells.Item($row, 13).Value = 0.00001
$ws.Cells.Item($row, 14).Value = 500
